# Refresh the Price (D) and Volume(1h) (E) columns for every row whose
# crypto snapshot values changed. Price values are stored as text in the
# sheet (matching the source data), so numeric-looking prices are entered
# with a leading apostrophe to force text entry and avoid Excel's
# automatic Text -> Number coercion (which would also drop trailing zeros,
# e.g. "4.80" -> 4.8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.932.36"
$ws.Range("D3").Value = "2.635.38"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'597.46"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "2.634.62"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'0.134"
$ws.Range("E10").Value = "  +9.52%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "'0.347"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'27.64"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").Value = "3.116.51"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "67.782.65"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "2.633.46"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "'11.45"
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").Value = "'372.08"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "'7.48"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "'4.80"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "'72.16"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'9.94"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").Value = "2.766.30"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "'579.14"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "'7.87"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "'157.90"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").Value = "'19.17"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("D41").Value = "'0.368"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "0.0₆0339"
$ws.Range("E43").Value = "  +16.98%  "
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("D45").Value = "'17.25"
$ws.Range("E45").Value = "  +5.73%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'40.26"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Value = "'156.09"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").Value = "'3.69"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").Value = "'21.96"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("D51").Value = "'1.70"
$ws.Range("E51").Value = "  -1.46%  "
